$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "DGEEC" source rows (23, 24) - they will be re-added
# further down the sheet (rows 29, 30) to make room for the new
# MSME classification table.
$ws.Range("A23").Clear()
$ws.Range("A24").Clear()

# New table header (row 18) - bold, matches the "title" style used
# elsewhere on the sheet (e.g. row 9).
$ws.Range("B18").Value = "Number of employees"
$ws.Range("B18").Font.Bold = $true
$ws.Range("C18").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("C18").Font.Bold = $true
$ws.Range("D18").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("D18").Font.Bold = $true

# Row 19 - Micro
$ws.Range("A19").Value = "Micro"
$ws.Range("B19").Value = "<5"
$ws.Range("C19").Value = "< G.23 Millionlon (USD4,400 approx.)"
$ws.Range("D19").Value = "< G.70 Millionlon (USD13,4450 approx.)"

# Row 20 - Small
$ws.Range("A20").Value = "Small"
$ws.Range("B20").Value = "6-20"
$ws.Range("C20").Value = "< G.92 Millionlon (USD17,700 approx.)"
$ws.Range("D20").Value = "< G.271 Millionlon (USD52,000 approx.)"

# Row 21 - Medium
$ws.Range("A21").Value = "Medium"
$ws.Range("B21").Value = "21-100"
$ws.Range("C21").Value = "< G.460 Millionlon (USD88,500 approx.)"
$ws.Range("D21").Value = "< G.1,355 Millionlon (USD 260,600approx.)"

# Row 22 - Large
$ws.Range("A22").Value = "Large"
$ws.Range("B22").Value = ">100"
$ws.Range("C22").Value = "> G.460 Millionlon (USD88,500 approx.)"
$ws.Range("D22").Value = "> G.1,355 Millionlon (USD 260,600approx.)"

# Re-add the DGEEC source attribution further down the sheet.
$ws.Range("A29").Value = "DGEEC"
$ws.Range("A29").Font.Bold = $true

$ws.Range("A30").Value = "Dirección General de Estadística, Encuestas y Censos (DGEEC), ""Censo Económico Nacional 2011"", 2013, p.57. Available at http://www.dgeec.gov.py/Publicaciones/Biblioteca/CEN2011/resultados_finales_CEN.pdf"
$ws.Range("A30").Font.Italic = $true
